# Update the cryptocurrency price/volume table to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'89.617.05"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "'3.031.28"

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'211.19"
$ws.Range("E5").Value = "  -2.01%  "

$ws.Range("D6").Value = "'612.35"
$ws.Range("E6").Value = "  -3.81%  "

$ws.Range("D7").Value = "'0.365"
$ws.Range("E7").Value = "  -8.01%  "

$ws.Range("D8").Value = "'0.895"
$ws.Range("E8").Value = "  +15.96%  "

$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").Value = "'3.028.65"
$ws.Range("E10").Value = "  -3.44%  "

$ws.Range("D11").Value = "'0.664"
$ws.Range("E11").Value = "  +18.44%  "

$ws.Range("E12").Value = "  +5.04%  "

$ws.Range("E13").Value = "  -4.99%  "

$ws.Range("E14").Value = "  -0.60%  "

$ws.Range("D15").Value = "'89.321.32"
$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("D16").Value = "'32.23"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").Value = "'3.582.06"
$ws.Range("E17").Value = "  -3.52%  "

$ws.Range("D18").Value = "'3.023.22"
$ws.Range("E18").Value = "  -3.80%  "

$ws.Range("D19").Value = "'3.32"
$ws.Range("E19").Value = "  -2.45%  "

$ws.Range("E20").Value = "  -4.26%  "

$ws.Range("D21").Value = "'13.38"
$ws.Range("E21").Value = "  +0.86%  "

$ws.Range("D22").Value = "'423.88"
$ws.Range("E22").Value = "  -0.52%  "

$ws.Range("D23").Value = "'8.26"
$ws.Range("E23").Value = "  -2.31%  "

$ws.Range("D24").Value = "'5.03"
$ws.Range("E24").Value = "  +2.15%  "

$ws.Range("D25").Value = "'5.33"
$ws.Range("E25").Value = "  -1.96%  "

$ws.Range("D26").Value = "'82.73"
$ws.Range("E26").Value = "  +0.84%  "

$ws.Range("D27").Value = "'11.54"
$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.22"
$ws.Range("E29").Value = "  +21.83%  "

$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").Value = "'0.161"
$ws.Range("E30").Value = "  +2.39%  "

$ws.Range("D31").Value = "'8.47"
$ws.Range("E31").Value = "  +3.41%  "

$ws.Range("D32").Value = "'3.73"
$ws.Range("E32").Value = "  -8.02%  "

$ws.Range("D33").Value = "'500.35"
$ws.Range("E33").Value = "  -1.31%  "

$ws.Range("D34").Value = "'6.62"
$ws.Range("E34").Value = "  -5.50%  "

$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").Value = "'1.81"
$ws.Range("E35").Value = "  -1.49%  "

$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'22.75"
$ws.Range("E36").Value = "  +2.60%  "

$ws.Range("E37").Value = "  -3.45%  "

$ws.Range("E38").Value = "  -9.47%  "

$ws.Range("D39").Value = "'22.26"
$ws.Range("E39").Value = "  -0.10%  "

$ws.Range("E40").Value = "  -0.18%  "

$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").Value = "'0.137"
$ws.Range("E42").Value = "  +4.66%  "

$ws.Range("D43").Value = "'0.357"
$ws.Range("E43").Value = "  -2.07%  "

$ws.Range("D44").Value = "'1.82"
$ws.Range("E44").Value = "  -2.93%  "

$ws.Range("D45").Value = "'142.95"
$ws.Range("E45").Value = "  -2.10%  "

$ws.Range("D46").Value = "'0.0694"
$ws.Range("E46").Value = "  +4.12%  "

$ws.Range("D47").Value = "'43.53"
$ws.Range("E47").Value = "  -0.33%  "

$ws.Range("E48").Value = "  +6.10%  "

$ws.Range("D49").Value = "'160.50"
$ws.Range("E49").Value = "  -2.55%  "

$ws.Range("D50").Value = "'1.22"
$ws.Range("E50").Value = "  +1.79%  "

$ws.Range("D51").Value = "'0.591"
$ws.Range("E51").Value = "  -1.49%  "

# Strip the quote-prefix formatting flag picked up above so the cells keep
# the workbook-default style (no explicit style index), matching the source.
$ws.Range("D2").ClearFormats()
$ws.Range("D3").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
